# Spankulator CPL update: version 5.3 of the board
# - Add new R16 component row (alphabetically inserted after R15)
# - Update reverse protection diode (D8/D9) positions
# - Update TP2 test point position (related to the diode change)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Spankulator-top-pos")

# --- Update D8 (row 30) and D9 (row 31) positions ---
$arrD8 = New-Object 'object[,]' 1,2
$arrD8[0,0] = 5.4
$arrD8[0,1] = 100.45
$ws.Range("B30:C30").Value = $arrD8

$ws.Range("C31").Value = 108.55

# --- Insert a new row for R16 between R15 (row 59) and R17 (old row 60) ---
$ws.Rows.Item(60).Insert()

$arrR16 = New-Object 'object[,]' 1,5
$arrR16[0,0] = "R16"
$arrR16[0,1] = 52.55
$arrR16[0,2] = 71.67
$arrR16[0,3] = 90
$arrR16[0,4] = "top"
$ws.Range("A60:E60").Value = $arrR16

# --- Update TP2 position (now shifted down to row 84) ---
$arrTP2 = New-Object 'object[,]' 1,2
$arrTP2[0,0] = 40.799999999999997
$arrTP2[0,1] = 103.8
$ws.Range("B84:C84").Value = $arrTP2

# --- Resize the query table to include the new row ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E94"))

# --- Update the ExternalData_1 defined name range to match the new extent ---
foreach ($n in $wb.Names) {
    if ($n.Name -like "*ExternalData_1*") {
        $n.RefersTo = "='Spankulator-top-pos'!`$A`$1:`$E`$94"
    }
}
